# Update the "今年日均" (this-year daily average) column H values on the
# "供应链放还款202310" (supply chain lending/repayment) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("供应链放还款202310")

$updates = @{
    "H2"  = 3935.85
    "H3"  = 894.83
    "H4"  = 6067.41
    "H5"  = 152.47
    "H6"  = 79.97
    "H7"  = 20398.65
    "H8"  = 76656.45
    "H9"  = 90.28
    "H10" = 17083.9
    "H11" = 18701.35
    "H12" = 19887.66
    "H13" = 240662.7
    "H14" = 1277.37
    "H15" = 62.84
    "H16" = 6772.51
    "H17" = 412724.24
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
